# Apply the "first block" calculations on the "Operating Conditions" sheet.
# This fills in the Io/Po/Qmax/Qavg style-section formulas (rows 17-22) that
# were previously left blank, which in turn feed the downstream cached
# formula results on the "Diode" sheet (D26/F26).
#
# Note: single-quoted strings are used for all formulas below so that
# PowerShell does not try to interpolate the "$" cell-reference anchors
# (e.g. $B$7) as variable references.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Operating Conditions")

# --- Row 17 : Io ---
$ws.Range("B17").Formula = '=B7*1.1'
$ws.Range("D17").Formula = '=B6/D16'
$ws.Range("F17").Formula = '=B6/F16'

# --- Row 18 : Po ---
$ws.Range("B18").Formula = '=B17*B16'
$ws.Range("D18").Formula = '=D17*D16'
$ws.Range("F18").Value = 15

# --- Row 19 : Qmax ---
$ws.Range("B19").Formula = '=B17^2 *B16'
$ws.Range("D19").Formula = '=D18^2/D16'
$ws.Range("F19").Formula = '=F18^2/F16'

# --- Row 20 : D ---
$ws.Range("B20").Value = 0
$ws.Range("D20").Formula = '=22.3'
$ws.Range("F20").Value = 53.38

# --- Row 21 : Qavg ---
$ws.Range("B21").Formula = '=B19/($B$7*$B$8)'

# --- Row 22 : Qavg ---
$ws.Range("B22").Formula = '=B19/($B$7*$B$8)'
$ws.Range("D22").Formula = '=D19/($B$7*$B$8)'
$ws.Range("F22").Formula = '=F19/($B$7*$B$8)'

# Update the saved selection to match the author's final cursor position.
[void]$ws.Activate()
[void]$ws.Range("D22").Select()

# Force a full recalculation so dependent sheets (e.g. Diode!D26/F26, which
# reference 'Operating Conditions'!D20 and F20) pick up the new cached values.
[void]$excel.CalculateFullRebuild()
